$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.9902682024985552
$ws.Range("F3").Value = 0.9524596519768238
$ws.Range("F4").Value = 0.9379133321344852
$ws.Range("F5").Value = 0.9486318528652191
$ws.Range("F6").Value = 0.9907087041065097
$ws.Range("F7").Value = 0.9975055702961981
$ws.Range("F8").Value = 0.9927204083651304
$ws.Range("F9").Value = 0.9908512523397803
$ws.Range("F10").Value = 0.9800759553909302
$ws.Range("F11").Value = 0.99901595688425
$ws.Range("F12").Value = 0.9577311463654041
$ws.Range("F13").Value = 0.9606174230575562
$ws.Range("F14").Value = 0.9980188477784395
$ws.Range("F15").Value = 0.9295283704996109
$ws.Range("F16").Value = 0.9929364589042962
$ws.Range("F17").Value = 0.6725634336471558
$ws.Range("F18").Value = 0.9918723702430725
$ws.Range("F19").Value = 0.9409340396523476
$ws.Range("F20").Value = 0.9906089622527361
$ws.Range("F21").Value = 0.9936257107183337
$ws.Range("F22").Value = 0.9977231787052006
$ws.Range("F23").Value = 0.9950504712760448
$ws.Range("F24").Value = 0.9988517602905631
$ws.Range("F25").Value = 0.9914289759472013
$ws.Range("F26").Value = 0.9459285624325275
$ws.Range("F27").Value = 0.9946284284815192
$ws.Range("F28").Value = 0.9773000385612249
$ws.Range("F29").Value = 0.9887322187423706
$ws.Range("F30").Value = 0.9919607043266296
$ws.Range("F31").Value = 0.9917110204696655
$ws.Range("F32").Value = 0.9973405003547668
$ws.Range("F33").Value = 0.9956170320510864
$ws.Range("F34").Value = 0.997776210308075
$ws.Range("F35").Value = 0.9933363199234009
$ws.Range("F36").Value = 0.9984061121940613
$ws.Range("F37").Value = 0.9940755367279053
$ws.Range("F38").Value = 0.9967204928398132
$ws.Range("F39").Value = 0.9871864318847656
$ws.Range("F40").Value = 0.9973031282424927
$ws.Range("F41").Value = 0.9910684823989868
$ws.Range("F42").Value = 0.9972319006919861
$ws.Range("F43").Value = 0.9948580265045166
$ws.Range("F44").Value = 0.996832549571991
$ws.Range("F45").Value = 0.9881040453910828
$ws.Range("F46").Value = 0.9783088564872742
$ws.Range("F47").Value = 0.9940879344940186
$ws.Range("F48").Value = 0.9962377548217773
$ws.Range("F49").Value = 0.9968886971473694
$ws.Range("F50").Value = 0.9887853860855103
$ws.Range("F51").Value = 0.9668812155723572
$ws.Range("F52").Value = 0.984312891960144
$ws.Range("F53").Value = 0.9976372718811035
$ws.Range("F54").Value = 0.9993595480918884
$ws.Range("F55").Value = 0.9958695769309998
$ws.Range("F56").Value = 0.9898996949195862
$ws.Range("F57").Value = 0.9982284903526306
$ws.Range("F58").Value = 0.996516227722168
$ws.Range("F59").Value = 0.9474180340766907
$ws.Range("F60").Value = 0.9950710535049438
$ws.Range("F61").Value = 0.978162944316864
$ws.Range("F62").Value = 0.9992918968200684
$ws.Range("F63").Value = 0.9973828196525574
$ws.Range("F64").Value = 0.9925822019577026
$ws.Range("F65").Value = 0.9919629096984863
$ws.Range("F66").Value = 0.9925715327262878
$ws.Range("F67").Value = 0.9977922439575195
$ws.Range("F68").Value = 0.9661133885383606
$ws.Range("F69").Value = 0.9857527613639832
$ws.Range("F70").Value = 0.9939227700233459
$ws.Range("F71").Value = 0.9989321827888489
$ws.Range("F72").Value = 0.9598271250724792
$ws.Range("F73").Value = 0.9753674864768982
$ws.Range("F74").Value = 0.9857174158096313
$ws.Range("F75").Value = 0.9984306693077087
$ws.Range("F76").Value = 0.9947405457496643
$ws.Range("F77").Value = 0.993002712726593
$ws.Range("F78").Value = 0.9848284125328064
$ws.Range("F79").Value = 0.9989218711853027
$ws.Range("F80").Value = 0.9986732006072998
$ws.Range("F81").Value = 0.9975624084472656
$ws.Range("F82").Value = 0.9972395896911621
$ws.Range("F83").Value = 0.9973642230033875
$ws.Range("F84").Value = 0.9867449402809143
$ws.Range("F85").Value = 0.9921119213104248
$ws.Range("F86").Value = 0.9978207349777222
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = $true
$ws.Range("F87").Value = 0.8810082376003265
$ws.Range("F88").Value = 0.9901557229459286
$ws.Range("F89").Value = 0.9989523759577423
$ws.Range("D90").Value = 1
$ws.Range("E90").Value = $false
$ws.Range("F90").Value = 0.9738199710845947
$ws.Range("F91").Value = 0.9438732266426086
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = $true
$ws.Range("F92").Value = 0.6744269132614136
$ws.Range("D93").Value = 1
$ws.Range("E93").Value = $false
$ws.Range("F93").Value = 0.8503650426864624
$ws.Range("F94").Value = 0.9237037301063538
$ws.Range("D95").Value = 1
$ws.Range("E95").Value = $false
$ws.Range("F95").Value = 0.9886132478713989
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = $true
$ws.Range("F96").Value = 0.5722082853317261
$ws.Range("F97").Value = 0.8486617207527161
$ws.Range("D98").Value = 1
$ws.Range("E98").Value = $false
$ws.Range("F98").Value = 0.6089860200881958
$ws.Range("F99").Value = 0.67630934715271
$ws.Range("F100").Value = 0.9806087203323841
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = $true
$ws.Range("F101").Value = 0.8542082756757736
$ws.Range("D102").Value = 1
$ws.Range("E102").Value = $false
$ws.Range("F102").Value = 0.9458279013633728
$ws.Range("F103").Value = 0.8669256567955017
$ws.Range("D104").Value = 0
$ws.Range("E104").Value = $true
$ws.Range("F104").Value = 0.9819917995482683
$ws.Range("F105").Value = 0.9098582714796066
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = $true
$ws.Range("F106").Value = 0.5959334671497345
$ws.Range("D107").Value = 1
$ws.Range("E107").Value = $false
$ws.Range("F107").Value = 0.9210028052330017
$ws.Range("F108").Value = 0.9449789710342884
$ws.Range("F109").Value = 0.8336544632911682
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = $true
$ws.Range("F110").Value = 0.9248353466391563
$ws.Range("F111").Value = 0.6510941982269287
$ws.Range("F112").Value = 0.7108362913131714
